$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0064018344450792568
$ws.Range("D2").Value = 0.062072178153578406
$ws.Range("E2").Value = 0.087167434194139748
$ws.Range("C3").Value = 0.0049528816253196131
$ws.Range("D3").Value = 0.04108042718816337
$ws.Range("E3").Value = 0.060495664608384965
$ws.Range("C4").Value = 0.0060960725181264707
$ws.Range("D4").Value = 0.11598512425920204
$ws.Range("E4").Value = 0.1398817904408467
$ws.Range("C5").Value = 0.0047255633040493789
$ws.Range("D5").Value = 0.073216881570390163
$ws.Range("E5").Value = 0.091741033858500826
$ws.Range("C6").Value = 0.0076090075985797805
$ws.Range("D6").Value = 0.16903851012283727
$ws.Range("E6").Value = 0.19886589718494874
$ws.Range("C7").Value = 0.0047428353121200971
$ws.Range("D7").Value = 0.11627335054979185
$ws.Range("E7").Value = 0.13486520890535678
$ws.Range("C8").Value = 0.0085249425335873277
$ws.Range("D8").Value = 0.20461849971486018
$ws.Range("E8").Value = 0.23803636102426756
$ws.Range("C9").Value = 0.0063052148017511567
$ws.Range("D9").Value = 0.15662869311725164
$ws.Range("E9").Value = 0.18134506060233047
$ws.Range("C10").Value = 0.010106448920241396
$ws.Range("D10").Value = 0.24665012120415147
$ws.Range("E10").Value = 0.28626750361073022
$ws.Range("C11").Value = 0.0079748456748932814
$ws.Range("D11").Value = 0.18110206609379653
$ws.Range("E11").Value = 0.21236336686386875
$ws.Range("C12").Value = 0.011083524908727353
$ws.Range("D12").Value = 0.2661494961961311
$ws.Range("E12").Value = 0.30959702640057846
$ws.Range("C13").Value = 0.0081490288485744616
$ws.Range("D13").Value = 0.20464844327170154
$ws.Range("E13").Value = 0.23659254002347854
$ws.Range("C14").Value = 0.010194617294331217
$ws.Range("D14").Value = 0.27139251137870607
$ws.Range("E14").Value = 0.31135551470713868
$ws.Range("C15").Value = 0.0088177616774583311
$ws.Range("D15").Value = 0.19321795874832454
$ws.Range("E15").Value = 0.22778348028382811
$ws.Range("C16").Value = 0.012985970856798304
$ws.Range("D16").Value = 0.20032620061788059
$ws.Range("E16").Value = 0.25123133825965643
$ws.Range("C17").Value = 0.010273851536478657
$ws.Range("D17").Value = 0.15449874300766686
$ws.Range("E17").Value = 0.19477211957720475
$ws.Range("C18").Value = 0.01521872914499676
$ws.Range("D18").Value = 0.058409727319396823
$ws.Range("E18").Value = 0.11806730012639226
$ws.Range("C19").Value = 0.012860236883300791
$ws.Range("D19").Value = 0.064899724121835034
$ws.Range("E19").Value = 0.1153117006756787
